$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal date-shaped string into a cell without letting
# Excel auto-convert it to a date serial. We stage the text as a formula
# result in a scratch cell, then Copy + PasteSpecial(values only) it into
# the destination so the destination keeps its original (unstyled) cell
# format, then we clean the scratch cell back up.
function Set-LiteralText($range, [string]$text) {
    $escaped = $text -replace '"', '""'
    $scratch = $ws.Range("ZZ1")
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $ws.Range($range).PasteSpecial(-4163)
    $scratch.ClearContents()
}

# Column map: A=Doi B=Date C=Title D=Abstract E=Authors F=ID G=ID Format
#             H=Date Accepted I=Other found locations J=Misc. Data

# --- Row 2: Fei Zhou et al. (Lancet) -------------------------------------
$ws.Range("E2").Value = "[Fei%Zhou%NULL%0, Ting%Yu%NULL%0, Ronghui%Du%NULL%0, Guohui%Fan%NULL%0, Ying%Liu%NULL%0, Zhibo%Liu%NULL%0, Jie%Xiang%NULL%0, Yeming%Wang%NULL%0, Bin%Song%NULL%0, Xiaoying%Gu%NULL%0, Lulu%Guan%NULL%0, Yuan%Wei%NULL%0, Hui%Li%NULL%0, Xudong%Wu%NULL%0, Jiuyang%Xu%NULL%0, Shengjin%Tu%NULL%0, Yi%Zhang%NULL%0, Hua%Chen%NULL%0, Bin%Cao%NULL%0]"
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = "Elsevier Ltd."

# --- Row 3: Litao Zhang et al. (D-dimer) ---------------------------------
$ws.Range("E3").Value = "[Litao%Zhang%NULL%0, Xinsheng%Yan%NULL%1, Qingkun%Fan%NULL%1, Haiyan%Liu%NULL%1, Xintian%Liu%NULL%1, Zejin%Liu%NULL%1, Zhenlu%Zhang%NULL%1]"
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = "International Society on Thrombosis and Haemostasis. Published by Elsevier Inc."

# --- Row 4: Rong-Hui Du et al. (Predictors of mortality) -----------------
$ws.Range("E4").Value = "[Rong-Hui%Du%NULL%0, Li-Rong%Liang%NULL%0, Cheng-Qing%Yang%NULL%0, Wen%Wang%NULL%0, Tan-Ze%Cao%NULL%0, Ming%Li%NULL%0, Guang-Yun%Guo%NULL%0, Juan%Du%NULL%0, Chun-Lan%Zheng%NULL%0, Qi%Zhu%NULL%0, Ming%Hu%NULL%0, Xu-Yan%Li%NULL%0, Peng%Peng%NULL%0, Huan-Zhong%Shi%NULL%0]"
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = "European Respiratory Society"

# --- Row 5: Jianlei Cao et al. (Clinical Features, 102 patients) ---------
$ws.Range("E5").Value = "[Jianlei%Cao%NULL%0, Wen-Jun%Tu%tuwenjun@irm-cams.ac.cn%0, Wenlin%Cheng%NULL%0, Lei%Yu%NULL%0, Ya-Kun%Liu%NULL%0, Xiaoyong%Hu%NULL%0, Qiang%Liu%NULL%0]"
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = "Oxford University Press"

# --- Row 6: Andrea Giacomelli et al. (30-day mortality) -------------------
$ws.Range("E6").Value = "[Andrea%Giacomelli%andrea.giacomelli@unimi.it%0, Anna Lisa%Ridolfo%NULL%0, Laura%Milazzo%NULL%0, Letizia%Oreni%NULL%0, Dario%Bernacchia%NULL%0, Matteo%Siano%NULL%0, Cecilia%Bonazzetti%NULL%0, Alice%Covizzi%NULL%0, Marco%Schiuma%NULL%0, Matteo%Passerini%NULL%0, Marco%Piscaglia%NULL%0, Massimo%Coen%NULL%0, Guido%Gubertini%NULL%0, Giuliano%Rizzardini%NULL%0, Chiara%Cogliati%NULL%0, Anna Maria%Brambilla%NULL%0, Riccardo%Colombo%NULL%0, Antonio%Castelli%NULL%0, Roberto%Rech%NULL%0, Agostino%Riva%NULL%0, Alessandro%Torre%NULL%0, Luca%Meroni%NULL%0, Stefano%Rusconi%NULL%0, Spinello%Antinori%NULL%0, Massimo%Galli%NULL%0]"
$ws.Range("I6").Value = ""
$ws.Range("J6").Value = "Elsevier Ltd."

# --- Row 7: Clinical Characteristics NYC (Paranjpe) -- reset to unknown --
$ws.Range("C7").Value = "Unknown Title"
$ws.Range("D7").Value = "Unknown Abstract"
$ws.Range("E7").Value = "[]"
$ws.Range("F7").Value = "not found"
$ws.Range("G7").Value = "N/A"
Set-LiteralText "H7" "1970-01-01"
$ws.Range("I7").Value = ""

# --- Row 8: Tao Chen et al. (113 deceased patients) -----------------------
$ws.Range("E8").Value = "[Tao%Chen%NULL%0, Di%Wu%NULL%0, Huilong%Chen%NULL%0, Weiming%Yan%NULL%0, Danlei%Yang%NULL%0, Guang%Chen%NULL%0, Ke%Ma%NULL%0, Dong%Xu%NULL%0, Haijing%Yu%NULL%0, Hongwu%Wang%NULL%0, Tao%Wang%NULL%0, Wei%Guo%NULL%0, Jia%Chen%NULL%0, Chen%Ding%NULL%0, Xiaoping%Zhang%NULL%0, Jiaquan%Huang%NULL%0, Meifang%Han%NULL%0, Shusheng%Li%NULL%0, Xiaoping%Luo%NULL%0, Jianping%Zhao%NULL%0, Qin%Ning%NULL%0]"
$ws.Range("I8").Value = ""
$ws.Range("J8").Value = "BMJ Publishing Group Ltd."

# --- Row 9: Characteristics of patients with COVID-19 (X. Luo) -- reset --
$ws.Range("C9").Value = "Unknown Title"
$ws.Range("D9").Value = "Unknown Abstract"
$ws.Range("E9").Value = "[]"
$ws.Range("F9").Value = "not found"
$ws.Range("G9").Value = "N/A"
Set-LiteralText "H9" "1970-01-01"
$ws.Range("I9").Value = ""

# --- Row 10: Retrospective study risk factors (Qingchun Yao) -- reset ----
$ws.Range("C10").Value = "Unknown Title"
$ws.Range("E10").Value = "[]"
$ws.Range("F10").Value = "not found"
$ws.Range("G10").Value = "N/A"
Set-LiteralText "H10" "1970-01-01"
$ws.Range("J10").Value = ""
